$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (Wire Nr, Product Nr) right
$ws.Range("A1").EntireColumn.Insert()

# Set the new header in A1
$ws.Range("A1").Value = "Kanban Nr"

# Update selection to match the target workbook state
$ws.Range("D10").Select()
